# PopulationByProvince_1986_2016.xlsx edit
# Commit message: "Addig dplyr relational examples"
#
# Observed changes:
#  1. Header cell E1 renamed from "Ulster (part of)" to "Ulster".
#  2. The province header row (B1:E1) is now center-aligned (was left-aligned).
#  3. The numeric data cells (B2:E8) are now center-aligned (were right-aligned).
#  4. The last user selection before saving was the whole of column F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the "Ulster (part of)" header to "Ulster"
$ws.Range("E1").Value = "Ulster"

# 2. Center the province header labels (Leinster, Munster, Connacht, Ulster)
$ws.Range("B1:E1").HorizontalAlignment = -4108   # xlCenter

# 3. Center the numeric population figures
$ws.Range("B2:E8").HorizontalAlignment = -4108   # xlCenter

# 4. Match the final selection state (column F was selected last)
$ws.Columns("F:F").Select()
